# The "JJ" initials textbox (author signature, leftover annotation) on the
# Login/Logic class-diagram slide is no longer wanted now that
# loginstorage.txt is documented to be created next to the jar file -
# remove it.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 1" -and $shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "JJ") {
        $shp.Delete()
    }
}
